# Insert a new weekly data row right after row 62 (new row becomes row 63),
# pushing the existing rows 63-104 down to 64-105.
# The new row duplicates row 62's data except for the date (column D),
# which gets a new weekly reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63..104 down by one to make room for the new row at 63.
$ws.Rows("63:63").Insert()

# Copy row 62's values into the freshly inserted row 63 (the row-insert
# already carried row 62's formatting down onto the new row).
for ($col = 1; $col -le 18; $col++) {
    $src = $ws.Cells.Item(62, $col)
    $dst = $ws.Cells.Item(63, $col)
    $dst.Value = $src.Value()
}

# New date for this weekly reading (serial 44452 = 2021-09-13).
$ws.Cells.Item(63, 4).Value = 44452
